$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($row, $a, $b, $c, $d, $vals)
    $ws.Range("A$row").Value = $a
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# Row 2: FAPs -> Agrp -> Mc4r -> ECs
Set-RowData 2 "FAPs" "Agrp" "Mc4r" "ECs" @(2, 0.6666666666666666, 0.7402036666666666, 2.220611, 0.1419214585000345, 0.1419214585000345, 1, 0.3333333333333333, 0.03919966666666667, 0.117599, 0.2591884862911046, 0.2591884862911046, 0.02901573699877778, 0.2611416329889999, 0.03678440800084976, 0.03678440800084977)

# Row 3: FAPs -> Agrp -> Mc4r -> MuSCs
Set-RowData 3 "FAPs" "Agrp" "Mc4r" "MuSCs" @(2, 0.6666666666666666, 0.7402036666666666, 2.220611, 0.1419214585000345, 0.1419214585000345, 2, 0.6666666666666666, 0.1120403333333333, 0.336121, 0.7408115137088953, 0.7408115137088953, 0.08293266554788889, 0.7463939899309999, 0.1051370504991847, 0.1051370504991847)

# Row 4: Inflammatory-Mac -> Agrp -> Mc4r -> ECs
Set-RowData 4 "Inflammatory-Mac" "Agrp" "Mc4r" "ECs" @(1, 0.3333333333333333, 0.3561633333333333, 1.06849, 0.06828825903893201, 0.06828825903893201, 1, 0.3333333333333333, 0.03919966666666667, 0.117599, 0.2591884862911046, 0.2591884862911046, 0.01396148394555556, 0.12565335551, 0.01769953049175563, 0.01769953049175563)

# Row 5: Inflammatory-Mac -> Agrp -> Mc4r -> MuSCs
Set-RowData 5 "Inflammatory-Mac" "Agrp" "Mc4r" "MuSCs" @(1, 0.3333333333333333, 0.3561633333333333, 1.06849, 0.06828825903893201, 0.06828825903893201, 2, 0.6666666666666666, 0.1120403333333333, 0.336121, 0.7408115137088953, 0.7408115137088953, 0.03990465858777778, 0.35914192729, 0.05058872854717638, 0.05058872854717638)

# Row 6: Neutrophils -> Agrp -> Mc4r -> ECs
Set-RowData 6 "Neutrophils" "Agrp" "Mc4r" "ECs" @(3, 1, 4.119219666666667, 12.357659, 0.7897902824610334, 0.7897902824610334, 1, 0.3333333333333333, 0.03919966666666667, 0.117599, 0.2591884862911046, 0.2591884862911046, 0.1614720378601111, 1.453248340741, 0.2047045477984992, 0.2047045477984992)

# Row 7: Neutrophils -> Agrp -> Mc4r -> MuSCs
Set-RowData 7 "Neutrophils" "Agrp" "Mc4r" "MuSCs" @(3, 1, 4.119219666666667, 12.357659, 0.7897902824610334, 0.7897902824610334, 2, 0.6666666666666666, 0.1120403333333333, 0.336121, 0.7408115137088953, 0.7408115137088953, 0.4615187445265556, 4.153668700739, 0.5850857346625342, 0.5850857346625342)
